# Generate Report for Handoff
# Updates the localization-status workbook:
#  - the 650f1eed... file moves to "In Translation"
#  - the a12e0c7a... file moves to "Ready for handoff"
#  - refreshed handoff/handback timestamps
#  - new "Error Detail" messages about stale handback files
#  - narrower Status column / wider Error Detail column

$wb = $excel.ActiveWorkbook

$statusInTranslation = "In Translation"
$statusReadyForHandoff = "Ready for handoff"

$newHandoffDateZhCn = "2016-10-13 13:55:21"
$newHandoffDateDeDe = "2016-10-13 13:55:33"

$errorDetail650f1eed = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af6da57df78594503dbd058d30d799a650731141/e2e/650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/870c1621b65d2f3fe583980256e85c6f97367043/e2e/650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.md."
$errorDetailA12e0c7a = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af6da57df78594503dbd058d30d799a650731141/e2e/a12e0c7a-e379-4e63-8710-ca2436491d04.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/870c1621b65d2f3fe583980256e85c6f97367043/e2e/a12e0c7a-e379-4e63-8710-ca2436491d04.md."

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusInTranslation
$wsOverview.Range("F2").Value = $statusInTranslation
$wsOverview.Range("G2").Value = $newHandoffDateDeDe

$wsOverview.Range("E3").Value = $statusReadyForHandoff
$wsOverview.Range("F3").Value = $statusReadyForHandoff
$wsOverview.Range("G3").Value = $newHandoffDateDeDe

$wsOverview.Columns.Item(5).ColumnWidth = 16.27
$wsOverview.Columns.Item(6).ColumnWidth = 16.27

# ---------- zh-cn sheet ----------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusInTranslation
$wsZhCn.Range("H2").Value = $newHandoffDateZhCn
$wsZhCn.Range("P2").Value = $errorDetail650f1eed

$wsZhCn.Range("C3").Value = $statusReadyForHandoff
$wsZhCn.Range("H3").Value = $newHandoffDateZhCn
$wsZhCn.Range("P3").Value = $errorDetailA12e0c7a

$wsZhCn.Columns.Item(3).ColumnWidth = 16.27
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------- de-de sheet ----------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusInTranslation
$wsDeDe.Range("H2").Value = $newHandoffDateDeDe
$wsDeDe.Range("P2").Value = $errorDetail650f1eed

$wsDeDe.Range("C3").Value = $statusReadyForHandoff
$wsDeDe.Range("H3").Value = $newHandoffDateDeDe
$wsDeDe.Range("P3").Value = $errorDetailA12e0c7a

$wsDeDe.Columns.Item(3).ColumnWidth = 16.27
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
